# Export and format DCF sensitivity table to Excel
$wb = $excel.ActiveWorkbook

$wsOutput = $wb.Worksheets.Item("output")

# ---------------------------------------------------------------------------
# 1. Add the new "sensitivity" worksheet as the last tab
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "sensitivity"

# ---------------------------------------------------------------------------
# 2. Column widths (row 1 indent column is narrow, data columns are wider)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 7
$ws.Columns.Item(2).ColumnWidth = 19
$ws.Columns.Item(3).ColumnWidth = 19
$ws.Columns.Item(4).ColumnWidth = 19

# ---------------------------------------------------------------------------
# 3. Terminal-growth headers (row 1, B1:D1) and WACC axis (A2:A6)
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = 0.02
$ws.Range("C1").Value = 0.025
$ws.Range("D1").Value = 0.03

$ws.Range("A2").Value = 0.08
$ws.Range("A3").Value = 0.085
$ws.Range("A4").Value = 0.09
$ws.Range("A5").Value = 0.095
$ws.Range("A6").Value = 0.1

# ---------------------------------------------------------------------------
# 4. Enterprise-value sensitivity grid (B2:D6)
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = 1419861275649.4
$ws.Range("C2").Value = 1522496183431.889
$ws.Range("D2").Value = 1645658072770.875

$ws.Range("B3").Value = 1310699804072.444
$ws.Range("C3").Value = 1395954929662.583
$ws.Range("D3").Value = 1496710987178.202

$ws.Range("B4").Value = 1217131849916.936
$ws.Range("C4").Value = 1288876028868.884
$ws.Range("D4").Value = 1372577570979.49

$ws.Range("B5").Value = 1136038731263.724
$ws.Range("C5").Value = 1197089098141.867
$ws.Range("D5").Value = 1267531829155.108

$ws.Range("B6").Value = 1065081436092.523
$ws.Range("C6").Value = 1117535853428.621
$ws.Range("D6").Value = 1177483758955.591

# ---------------------------------------------------------------------------
# 5. Formatting: bold / bordered / centered header row + axis column
#    (mirrors the "Metric"/"Value" header style already used on "output")
# ---------------------------------------------------------------------------
$headerRow = $ws.Range("B1:D1")
$headerRow.Font.Bold = $true
$headerRow.Borders.LineStyle = 1
$headerRow.HorizontalAlignment = -4108
$headerRow.VerticalAlignment = -4160

$headerCol = $ws.Range("A2:A6")
$headerCol.Font.Bold = $true
$headerCol.Borders.LineStyle = 1
$headerCol.HorizontalAlignment = -4108
$headerCol.VerticalAlignment = -4160

# Dollar-formatted value grid (same currency format as output!B2)
$valueRange = $ws.Range("B2:D6")
$valueRange.NumberFormat = "$#,##0"

# ---------------------------------------------------------------------------
# 6. Conditional formatting - 3 colour scale (red/yellow/green) on the grid
# ---------------------------------------------------------------------------
$cs = $valueRange.FormatConditions.AddColorScale(3)
$cs.ColorScaleCriteria.Item(1).Type = 1
$cs.ColorScaleCriteria.Item(2).Type = 4
$cs.ColorScaleCriteria.Item(2).Value = 50
$cs.ColorScaleCriteria.Item(3).Type = 2

# ---------------------------------------------------------------------------
# 7. Make the new sheet the active tab (matches author's saved view state)
# ---------------------------------------------------------------------------
$ws.Range("A1").Select() | Out-Null
$ws.Activate() | Out-Null
